$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" on every sheet ---
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- Column width changes (status columns got narrower once the longer
#     "Ready for handoff" text was replaced by the shorter "In Translation") ---
# NB: the target OOXML width is 13.4101845877511 "characters"; this runtime
# quantizes ColumnWidth writes to 1/6-character steps (Width = round(x*6)/6 + 5/6),
# so 12.5 is the input that lands on the nearest achievable grid point (13.3333...).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").ColumnWidth = 12.5
